# Edit script: update EC database - replace worker detail rows with new data
# (Roberto Carlos Vergara Buelvas total + Efrain/Olmar/Diego x 4 periods),
# update header totals, and remove the now-unused extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the bold "last row" formatting (currently on row 33) onto row 28,
#        which will become the new bold last data row. ---
$ws.Range("B33:J33").Copy() | Out-Null
$ws.Range("B28:J28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Update header summary values ---
$ws.Range("E11").Value = 1110568    # VALOR MORA
$ws.Range("C13").Value = 4          # Cant. Trabajadores

# --- 3. Overwrite the worker detail rows (16-28) with the new data set ---

# Row 16: Roberto Carlos Vergara Buelvas / period 2411 (kept, not bold)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "92556796"
$ws.Range("D16").Value = "ROBERTO CARLOS VERGARA BUELVAS"
$ws.Range("E16").Value = "2411"
$ws.Range("F16").Value = 139672
$ws.Range("G16").Value = 3673356

# Rows 17-28: Efrain / Olmar / Diego across periods 2505, 2506, 2507, 2508
$workers = @(
    @("1050549994", "EFRAIN ANTONIO LOPEZ HERNANDEZ"),
    @("1050542558", "OLMAR ALCOCER GONZALEZ"),
    @("1050553847", "DIEGO FELIPE MILLAN SALGADO")
)
$periods = @("2505", "2506", "2507", "2508")

$row = 17
foreach ($period in $periods) {
    foreach ($worker in $workers) {
        $ws.Range("B$row").Value = "CC"
        $ws.Range("C$row").Value = $worker[0]
        $ws.Range("D$row").Value = $worker[1]
        $ws.Range("E$row").Value = $period
        $ws.Range("F$row").Value = 80908
        $ws.Range("G$row").Value = 2022676
        $row = $row + 1
    }
}

# --- 4. Remove the now unused rows 29-33 (old extra workers + old totals row),
#        which shifts the signature block (previously rows 38-39) up to rows 33-34 ---
$ws.Range("B29:J33").EntireRow.Delete() | Out-Null
